# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 15
$ws1.Range("F3").Value  = 173
$ws1.Range("F4").Value  = 83
$ws1.Range("F6").Value  = 546
$ws1.Range("F7").Value  = 1682
$ws1.Range("F8").Value  = 22
$ws1.Range("F10").Value = 29
$ws1.Range("F11").Value = 1647
$ws1.Range("F13").Value = 70
$ws1.Range("F14").Value = 405
$ws1.Range("F18").Value = 24
$ws1.Range("F20").Value = 51
$ws1.Range("F21").Value = 240
$ws1.Range("F23").Value = 160
$ws1.Range("F24").Value = 226
$ws1.Range("F25").Value = 237

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 15
$ws4.Range("F3").Value  = 173
$ws4.Range("F4").Value  = 83
$ws4.Range("F6").Value  = 546
$ws4.Range("F7").Value  = 1682
$ws4.Range("F9").Value  = 22
$ws4.Range("F11").Value = 29
$ws4.Range("F12").Value = 1647
$ws4.Range("F14").Value = 70
$ws4.Range("F15").Value = 405
$ws4.Range("F19").Value = 24
$ws4.Range("F21").Value = 51
$ws4.Range("F22").Value = 240
$ws4.Range("F24").Value = 160
$ws4.Range("F25").Value = 226
$ws4.Range("F26").Value = 237
